$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.064.44'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '2.221.47'
$ws.Range('E3').Value = '  -1.22%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('E5').Value = '  -1.79%  '
$ws.Range('D6').Value = '''0.628'
$ws.Range('E6').Value = '  +1.10%  '
$ws.Range('E7').Value = '  -0.50%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').Value = '''0.616'
$ws.Range('E9').Value = '  -0.79%  '
$ws.Range('E10').Value = '  +4.57%  '
$ws.Range('E11').Value = '  +1.95%  '
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').Value = '''0.104'
$ws.Range('E13').Value = '  +1.17%  '
$ws.Range('D14').Value = '2.551.46'
$ws.Range('E14').Value = '  -1.22%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '''14.33'
$ws.Range('E15').Value = '  -1.14%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = '''0.849'
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('D17').Value = '2.246.06'
$ws.Range('E17').Value = '  -0.21%  '
$ws.Range('D18').Value = '41.943.11'
$ws.Range('E18').Value = '  -0.44%  '
$ws.Range('E19').Value = '  +12.25%  '
$ws.Range('D20').Value = '''6.18'
$ws.Range('E20').Value = '  +0.76%  '
$ws.Range('D21').Value = '''72.68'
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('D22').Value = '''10.74'
$ws.Range('E22').Value = '  +20.58%  '
$ws.Range('D23').Value = '''230.33'
$ws.Range('E23').Value = '  -0.76%  '
$ws.Range('E24').Value = '  -6.61%  '
$ws.Range('D25').Value = '''11.63'
$ws.Range('E25').Value = '  +3.62%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  -1.86%  '
$ws.Range('D28').Value = '''2.28'
$ws.Range('E28').Value = '  -1.31%  '
$ws.Range('E29').Value = '  +2.49%  '
$ws.Range('D30').Value = '''167.05'
$ws.Range('E30').Value = '  -1.50%  '
$ws.Range('E31').Value = '  -0.15%  '
$ws.Range('D32').Value = '''5.72'
$ws.Range('E32').Value = '  +10.62%  '
$ws.Range('D33').Value = '''0.0798'
$ws.Range('E33').Value = '  -3.49%  '
$ws.Range('E34').Value = '  +0.90%  '
$ws.Range('E35').Value = '  -4.25%  '
$ws.Range('D36').Value = '''29.12'
$ws.Range('E36').Value = '  -4.04%  '
$ws.Range('D37').Value = '''4.31'
$ws.Range('E37').Value = '  -5.29%  '
$ws.Range('D38').Value = '''0.0304'
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('D39').Value = '''12.95'
$ws.Range('E39').Value = '  -3.96%  '
$ws.Range('D40').Value = '''67.09'
$ws.Range('E40').Value = '  +8.72%  '
$ws.Range('D41').Value = '''2.13'
$ws.Range('E41').Value = '  -2.62%  '
$ws.Range('E42').Value = '  -2.38%  '
$ws.Range('D43').Value = '''0.203'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').Value = '''8.82'
$ws.Range('E44').Value = '  +2.16%  '
$ws.Range('D45').Value = '''104.77'
$ws.Range('E45').Value = '  -3.91%  '
$ws.Range('E46').Value = '  +0.06%  '
$ws.Range('E47').Value = '  +5.55%  '
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('E50').Value = '  +0.21%  '
$ws.Range('D51').Value = '2.430.14'
$ws.Range('E51').Value = '  -1.13%  '
